$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 7-11 which are no longer part of the data (shifts rows up, nothing below anyway)
$ws.Rows("7:11").Delete()

# Update rows 2-6 with the refreshed group/count values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 129

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 113

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 81

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 71

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 62
